# Update gh-pages output data (refreshed crawl numbers / cover image)
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId=1) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 3122
$wsExhibit.Range("F4").Value = 147
$wsExhibit.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202410/ccEfc1521728888008037.jpeg"
$wsExhibit.Range("F5").Value = 115

# --- Sheet "全部类型" (sheetId=4) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 3122
$wsAll.Range("F8").Value = 147
$wsAll.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202410/ccEfc1521728888008037.jpeg"
$wsAll.Range("F10").Value = 115
